# Update cryptocurrency price/volume data per the latest refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.711.84"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "1.645.89"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Formula = "'213.30"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Formula = "'23.05"
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.879.44"
$ws.Range("D13").Value = "1.645.91"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Formula = "'64.23"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "27.684.84"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Formula = "'230.27"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Formula = "'7.66"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Formula = "'10.01"
$ws.Range("E23").Value = "  +6.90%  "
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").Formula = "'149.02"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Formula = "'15.66"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Formula = "'1.00"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Formula = "'3.31"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "1.440.63"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Formula = "'0.572"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Formula = "'0.884"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Formula = "'0.904"
$ws.Range("E40").Value = "  +14.87%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D43").Formula = "'5.66"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").Formula = "'65.60"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").Value = "1.788.60"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Formula = "'86.35"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("D50").Formula = "'0.0989"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").Formula = "'7.73"
$ws.Range("E51").Value = "  +0.16%  "
